# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change: cell B11 on the "Rules" sheet held the text "R40"; it now holds
# the text "1" (still a plain text value, not a number — a new shared
# string is appended for it). The cell's existing style/formatting must be
# left untouched.
#
# Assigning Range.Value = "1" directly would store a NUMBER (Excel infers
# type from the literal), which also happens to reallocate the cell's style
# record. To keep both the text type and the original style intact we:
#   1. stash a copy of B11's current formatting in an unused scratch cell,
#   2. write the new value with a leading apostrophe so Excel stores it as
#      literal text ("1") instead of the number 1,
#   3. paste the stashed formatting back onto B11 (formats only, so the
#      freshly-written text value is preserved),
#   4. clean up the scratch cell so it leaves no trace in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target  = $ws.Range("B11")
$scratch = $ws.Range("ZZ1")

$target.Copy($scratch)          # remember B11's current formatting
$target.Value = "'1"            # force a text value of "1" (not the number 1)

$scratch.Copy()
$target.PasteSpecial(-4122)     # xlPasteFormats: restore the original look
$scratch.Clear()                # leave no trace of the scratch cell

$excel.CutCopyMode = $false
